$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures to the latest scrape. Column D values are written with a
# leading apostrophe (quote-prefix) so Excel keeps them as literal text
# -- otherwise numeric-looking strings like "1.00" would be silently
# re-parsed as numbers and lose their trailing zeros/formatting. The
# cell style is reset to "Normal" right after so no stray text-format
# style gets left behind on the cell.

$ws.Range("D2").Value = '''53.407.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.97%  '

$ws.Range("D3").Value = '''3.152.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''398.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.69%  '

$ws.Range("D6").Value = '''107.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.14%  '

$ws.Range("D7").Value = '''0.546'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = '''0.610'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.87%  '

$ws.Range("D10").Value = '''38.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.16%  '

$ws.Range("D11").Value = '''0.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("E12").Value = '  +1.25%  '

$ws.Range("D13").Value = '''3.647.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.41%  '

$ws.Range("D14").Value = '''18.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.37%  '

$ws.Range("D15").Value = '''7.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.19%  '

$ws.Range("E16").Value = '  +8.98%  '

$ws.Range("D17").Value = '''3.150.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.28%  '

$ws.Range("D18").Value = '''10.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").Value = '''53.290.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("D20").Value = '''3.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.76%  '

$ws.Range("D21").Value = '''12.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.84%  '

$ws.Range("E22").Value = '  +1.36%  '

$ws.Range("D23").Value = '''71.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").Value = '''271.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.22%  '

$ws.Range("E25").Value = '  +2.66%  '

$ws.Range("D26").Value = '''8.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.75%  '

$ws.Range("D27").Value = '''27.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.42%  '

$ws.Range("D28").Value = '''7.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.79%  '

$ws.Range("D29").Value = '''0.170'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("E31").Value = '  +2.60%  '

$ws.Range("D32").Value = '''11.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.33%  '

$ws.Range("D33").Value = '''37.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.73%  '

$ws.Range("D34").Value = '''0.0493'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.73%  '

$ws.Range("D35").Value = '''2.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.69%  '

$ws.Range("D36").Value = '''50.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").Value = '''3.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.40%  '

$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("D40").Value = '''4.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.70%  '

$ws.Range("E41").Value = '  +0.97%  '

$ws.Range("D42").Value = '''17.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.12%  '

$ws.Range("D43").Value = '''1.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.37%  '

$ws.Range("D44").Value = '''130.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.44%  '

$ws.Range("E45").Value = '  +1.50%  '

$ws.Range("D46").Value = '''22.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.43%  '

$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("D48").Value = '''2.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").Value = '''2.093.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.28%  '

$ws.Range("D50").Value = '''0.0510'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +27.71%  '

$ws.Range("D51").Value = '''0.0334'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.79%  '
